$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S: header year 2022 in S4 (same style as R4/Q4)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Update existing values in row 5
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.899999999999999

# New value in S5 (same style as R5)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 13.5

# Update selection to match new active range
$ws.Range("S7:S8").Select()
